# [FEATURE] Add settings screen and wip notifications and friend requests
#
# This script reproduces the edit made to Artefakte/Zeiterfassung/Arbeitsmatrix.xlsx:
#   - A new "Einstellungen Screen, Icons und Ideensammlung" task entry is added
#     to the Arbeitsmatrix sheet (what used to be a blank spacer row, row 129,
#     becomes a filled-in data row, mirroring the rows above it).
#   - Four new blank spacer rows are inserted below it (so the totals block that
#     used to start at row 130 now starts at row 134).
#   - The Prefix data-validation range is extended to cover the new blank rows.
#   - The sheet view scrolls down / selection moves to reflect the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# ---------------------------------------------------------------------------
# 1) Make room: insert 4 new blank rows just above the "totals" block.
#    Old row 130 (Stunden insgesamt ...) becomes row 134, etc.
# ---------------------------------------------------------------------------
$ws.Rows("130:133").Insert()

# ---------------------------------------------------------------------------
# 2) Turn the (now former) blank spacer row 129 into a real data row, mirroring
#    the formatting of the row above it (row 128).
# ---------------------------------------------------------------------------
$ws.Range("A128:K128").Copy()
$ws.Range("A129:K129").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H129").ClearContents()

$ws.Range("A129").Value = 22
$ws.Range("B129").Value = "Interface Design"
$ws.Range("C129").Value = "MockUps"
$ws.Range("D129").Value = "[FEATURE]"
$ws.Range("E129").Value = "Einstellungen Screen, Icons und Ideensammlung"
$ws.Range("F129").Value = "09/20/2021"
$ws.Range("G129").Value = "10/12/2021"
$ws.Range("I129").Formula = "=ROUNDUP(((SUM(K129-J129)*24*60/60)/0.25),0)*0.25"
$ws.Range("J129").Value = 0.75
$ws.Range("K129").Value = 0.82291666666666663

# ---------------------------------------------------------------------------
# 3) Extend the Prefix data-validation range so the new blank rows (130:133)
#    keep the same dropdown validation the old blank rows had (previously
#    D115:D129, now D115:D133).
# ---------------------------------------------------------------------------
$dvRange = $ws.Range("D115:D133")
$dvRange.Validation.Delete()
$dvRange.Validation.Add(3, 1, 1, "=`$N`$3:`$N`$6")
$dvRange.Validation.InCellDropdown = $true
$dvRange.Validation.ErrorTitle = "Prefix nicht unterstützt"
$dvRange.Validation.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden_x000a_"
$dvRange.Validation.PromptTitle = "Prefix"
$dvRange.Validation.InputMessage = "Wählen Sie einen Prefix aus"

# ---------------------------------------------------------------------------
# 4) Update the view: scroll down towards the new rows and move the selection.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 116
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E131").Select() | Out-Null
